# Orders template: turn row 2 into a filled-in example order, add three
# blank rows (3-5) pre-populated with the same boilerplate numbers/times
# for user data entry, and retune several column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column width adjustments.
# Excel's Range/Columns.ColumnWidth property is expressed in "characters
# of the Normal style font" while the raw OOXML <col width="..."> stores
# ColumnWidth + 5/6 (the standard internal column-padding constant for
# the default Calibri 11 font). Subtract 5/6 so the saved file lands
# exactly on the desired OOXML width.
# ---------------------------------------------------------------------
$pad = 5 / 6

function Set-ColWidth($colIndex, $targetOoxmlWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetOoxmlWidth - $pad
}

Set-ColWidth 1 12    # A: 6  -> 12
Set-ColWidth 3 5     # C: 16 -> 5
Set-ColWidth 4 11    # D: 9  -> 11
Set-ColWidth 5 11    # E: 9  -> 11
Set-ColWidth 9 6     # I: 5  -> 6
Set-ColWidth 10 10   # J: 6  -> 10
Set-ColWidth 17 7    # Q: 8  -> 7
Set-ColWidth 18 6    # R: 8  -> 6
Set-ColWidth 19 8    # S: 6  -> 8

# ---------------------------------------------------------------------
# Helper: write literal text into a cell without Excel re-interpreting
# look-alike dates/numbers (e.g. "2026-01-19" or "09:00") -- a leading
# apostrophe forces text, then resetting the style back to "Normal"
# drops the transient quote-prefix formatting Excel applies so the cell
# ends up indistinguishable from a plain text cell typed via the OOXML.
# An empty string ("'" with nothing after it) likewise yields a real,
# present-but-empty text cell instead of Excel silently clearing it.
# ---------------------------------------------------------------------
function Set-Text($ref, $text) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($ref).Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 2: fill in a concrete worked example (was mostly blank / generic
# placeholder instructions).
# ---------------------------------------------------------------------
Set-Text "A2" "예시-ORD-001"
Set-Text "B2" "2026-01-19"
Set-Text "C2" "냉동"
Set-Text "D2" "CUST-0001"
Set-Text "E2" "CUST-0002"
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 5000
$ws.Range("H2").Value = 15
Set-Text "I2" "냉동식품"
Set-Text "J2" "PROD-001"
Set-Text "K2" "09:00"
Set-Text "L2" "12:00"
Set-Text "M2" "14:00"
Set-Text "N2" "17:00"
Set-Text "O2" "2026-01-19"
$ws.Range("P2").Value = 5
Set-Text "Q2" "Y"
Set-Text "R2" "Y"
Set-Text "S2" "깨지기 쉬움"

# ---------------------------------------------------------------------
# Rows 3-5: three blank rows for the user's own data entry. Only the
# boilerplate defaults (dates, quantities, standard time windows, Y/Y
# flags) are pre-filled; the identifying/free-text columns (A, C, D, E,
# I, J, S) are left as empty text cells for the user to complete.
# ---------------------------------------------------------------------
foreach ($r in 3..5) {
    Set-Text "A$r" ""
    Set-Text "B$r" "2026-01-19"
    Set-Text "C$r" ""
    Set-Text "D$r" ""
    Set-Text "E$r" ""
    $ws.Range("F$r").Value = 10
    $ws.Range("G$r").Value = 5000
    $ws.Range("H$r").Value = 15
    Set-Text "I$r" ""
    Set-Text "J$r" ""
    Set-Text "K$r" "09:00"
    Set-Text "L$r" "12:00"
    Set-Text "M$r" "14:00"
    Set-Text "N$r" "17:00"
    Set-Text "O$r" "2026-01-19"
    $ws.Range("P$r").Value = 5
    Set-Text "Q$r" "Y"
    Set-Text "R$r" "Y"
    Set-Text "S$r" ""
}
